$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price / volume data to reflect the latest scrape.
# Column D ("Price") values that look numeric are entered with a leading
# apostrophe so Excel stores them as text (matching the sheet's existing
# convention of keeping prices as literal strings, e.g. "68.313.39").

$ws.Range("D2").Value = "68.313.39"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "3.564.24"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'619.47"
$ws.Range("E5").Value = "  +2.27%  "
$ws.Range("D6").Value = "'155.32"
$ws.Range("E6").Value = "  +4.18%  "
$ws.Range("D7").Value = "3.561.59"
$ws.Range("E7").Value = "  +1.81%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("E10").Value = "  +5.34%  "
$ws.Range("E11").Value = "  +6.79%  "
$ws.Range("E12").Value = "  +4.02%  "
$ws.Range("D13").Value = "'33.29"
$ws.Range("E13").Value = "  +5.87%  "
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "4.167.40"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "3.567.64"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").Value = "68.352.07"
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("E19").Value = "  +5.87%  "
$ws.Range("D20").Value = "'16.03"
$ws.Range("E20").Value = "  +7.13%  "
$ws.Range("D21").Value = "'10.03"
$ws.Range("E21").Value = "  +12.04%  "
$ws.Range("D22").Value = "'454.09"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").Value = "'0.643"
$ws.Range("E23").Value = "  +4.09%  "
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("E25").Value = "  +2.35%  "
$ws.Range("D26").Value = "3.707.85"
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").Value = "'9.21"
$ws.Range("E28").Value = "  +11.85%  "
$ws.Range("D29").Value = "'10.55"
$ws.Range("E29").Value = "  +3.71%  "
$ws.Range("E30").Value = "  +10.79%  "
$ws.Range("E31").Value = "  +3.40%  "
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("D33").Value = "'0.998"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "'6.38"
$ws.Range("E34").Value = "  +4.13%  "
$ws.Range("D35").Value = "'26.16"
$ws.Range("E35").Value = "  +1.98%  "
$ws.Range("E36").Value = "  +3.84%  "
$ws.Range("D37").Value = "3.558.26"
$ws.Range("E37").Value = "  +1.99%  "
$ws.Range("E38").Value = "  +3.37%  "
$ws.Range("E39").Value = "  +9.18%  "
$ws.Range("D41").Value = "'181.29"
$ws.Range("E41").Value = "  +4.56%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.0917"
$ws.Range("E42").Value = "  +4.94%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "'5.60"
$ws.Range("E44").Value = "  +3.60%  "
$ws.Range("D45").Value = "'31.17"
$ws.Range("E45").Value = "  +12.94%  "
$ws.Range("D46").Value = "'0.899"
$ws.Range("D47").Value = "'46.22"
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("D48").Value = "'1.34"
$ws.Range("E48").Value = "  +4.29%  "
$ws.Range("D49").Value = "'2.67"
$ws.Range("E49").Value = "  +4.51%  "
$ws.Range("E50").Value = "  +3.56%  "
$ws.Range("D51").Value = "'0.263"
$ws.Range("E51").Value = "  +7.94%  "

Write-Host "Cryptos list updated"
